$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A50").Value = "ХЕЧ ТЕК"
$ws.Range("B50").Value = "ХЕЧ ТЕК Україна"
